$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Resultados" intro sentence: mention the single-process run too.
#    " com 4 processos e 4 threads respectivamente"
#      -> " com um único processo, com 4 processos e com 4 threads
#           atuando, respectivamente"
# ---------------------------------------------------------------------
[void]$d.Content.Find.Execute(
    " com 4 processos e 4 threads respectivamente",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " com um único processo, com 4 processos e com 4 threads atuando, respectivamente",
    2
)

# ---------------------------------------------------------------------
# 2. Add the new single-process measurement as the first bullet of the
#    results list (before the existing "3.0094 ± 0.1697 [s]" line).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "3.0094*") {
        $para.Range.InsertParagraphBefore()
        $newItem = $d.Paragraphs.Item($i)
        $newItem.Range.Text = "6.6762 ± 0.1577 [s]"
        break
    }
}

# ---------------------------------------------------------------------
# 3. Discussion paragraph: add a closing remark about the speed-up
#    gained from parallelism.
# ---------------------------------------------------------------------
[void]$d.Content.Find.Execute(
    "considerando os desvios padrão, vemos que há uma faixa de valores onde os processos foram mais rápidos que os threads.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "considerando os desvios padrão, vemos que há uma faixa de valores onde os processos foram mais rápidos que os threads. Além disso, fica claro que a aplicação do paralelismo nessa tarefa obteve um grande aumento de desempenho.",
    2
)

# ---------------------------------------------------------------------
# 4. Join the two sentences about thread efficiency into one:
#    "...próximos. E que a aplicação..." -> "...próximos, e que a aplicação..."
# ---------------------------------------------------------------------
[void]$d.Content.Find.Execute(
    "próximos. E que a aplicação que utilizou threads obteve uma eficiência maior,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "próximos, e que a aplicação que utilizou threads obteve uma eficiência maior,",
    2
)
